$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2: balanced-pie-with-filler example data ---

# Header row (right-aligned labels)
$ws2.Range("C3").Value = "A"
$ws2.Range("D3").Value = "B"
$ws2.Range("E3").Value = "Filler"
$ws2.Range("C3:E3").HorizontalAlignment = -4152   # xlRight

# First example: Before / Proposal / After
$ws2.Range("B4").Value = "Before"
$ws2.Range("C4").Value = 80
$ws2.Range("D4").Value = 15
$ws2.Range("E4").Value = 5
$ws2.Range("F4").Formula = "=SUM(C4:E4)"

$ws2.Range("B5").Value = "Proposal"
$ws2.Range("D5").Value = 30

$ws2.Range("B6").Value = "After"
$ws2.Range("C6").Value = 80
$ws2.Range("D6").Value = 20
$ws2.Range("E6").Value = 0
$ws2.Range("F6").Formula = "=SUM(C6:E6)"

# Second example: Before / Proposal / After
$ws2.Range("B8").Value = "Before"
$ws2.Range("C8").Value = 60
$ws2.Range("D8").Value = 20
$ws2.Range("E8").Value = 20
$ws2.Range("F8").Formula = "=SUM(C8:E8)"

$ws2.Range("B9").Value = "Proposal"
$ws2.Range("D9").Value = 25

$ws2.Range("B10").Value = "After"
$ws2.Range("C10").Value = 60
$ws2.Range("D10").Value = 25
$ws2.Range("E10").Value = 15

# Filler ratio helper cells
$ws2.Range("L17").Value = 40
$ws2.Range("M17").Formula = "=40/90"

$ws2.Range("L18").Value = 50
$ws2.Range("M18").Formula = "=50/90"

# --- View / selection state ---
# Sheet1 keeps its selection but is no longer the tab in view; scroll its
# window so B43 is the top-left visible cell (best-effort - mirrors the
# author scrolling down to the new example before switching tabs).
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 2
$ws1.Range("E60").Select()

# Sheet2 becomes the active/selected tab, scrolled/selected at M19.
$ws2.Activate()
$ws2.Range("M19").Select()
